$d = $word.ActiveDocument

# Remove the paragraphs: "Hey loko rescatate...", the blank paragraph after it,
# and the "xD XD" paragraph. These are paragraphs 2, 3 and 4 of the document.
$start = $d.Paragraphs(2).Range.Start
$end = $d.Paragraphs(4).Range.End
$d.Range($start, $end).Delete()
